$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# A3: update task description text (HTML + JavaScript -> HTML + CSS)
$ws.Range("A3").Value = "Criação e integração da interface com HTML e CSS"

# Row 10 (T9 - Realizar testes de Integração): fill in time spent / total
$ws.Range("F10").Value = "1h+30min"
$ws.Range("G10").Value = "1h30"

# Row 11 (T10 - Adaptações, se necessárias): update people, time spent, total
$ws.Range("E11").Value = "Ana, Laís e Beatriz <3"
$ws.Range("F11").Value = "1h+1h30"
$ws.Range("G11").Value = "2h30"

# Update the active selection to A3
$ws.Range("A3").Select()
